# Scheduled market-data refresh: push updated FFXIV Market Board price/profit
# snapshots (currentAveragePrice*, LevePrice*, LeveProfit*) into each job sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 10020.223
$ws.Range("I34").Value = 10020.223
$ws.Range("K34").Value = 10020.223
$ws.Range("M34").Value = -9817.223
$ws.Range("H36").Value = 10020.223
$ws.Range("I36").Value = 10020.223
$ws.Range("K36").Value = 10020.223
$ws.Range("M36").Value = -9305.223
$ws.Range("H70").Value = 3254.7144
$ws.Range("J70").Value = 3636.6
$ws.Range("L70").Value = 10909.8
$ws.Range("N70").Value = -11449.8
$ws.Range("H73").Value = 3254.7144
$ws.Range("J73").Value = 3636.6
$ws.Range("L73").Value = 10909.8
$ws.Range("N73").Value = -12781.8
$ws.Range("H92").Value = 401
$ws.Range("I92").Value = 396.18182
$ws.Range("J92").Value = 427.5
$ws.Range("K92").Value = 396.18182
$ws.Range("L92").Value = 427.5
$ws.Range("M92").Value = 851.81818
$ws.Range("N92").Value = -2923.5
$ws.Range("H100").Value = 8594.5
$ws.Range("I100").Value = 4228.75
$ws.Range("J100").Value = 11505
$ws.Range("K100").Value = 4228.75
$ws.Range("L100").Value = 11505
$ws.Range("M100").Value = -3687.75
$ws.Range("N100").Value = -12587
$ws.Range("H137").Value = 1347.38
$ws.Range("I137").Value = 1870.2963
$ws.Range("J137").Value = 1153.9727
$ws.Range("K137").Value = 5610.8889
$ws.Range("L137").Value = 3461.9181
$ws.Range("M137").Value = -3060.8889
$ws.Range("N137").Value = -8561.918099999999
$ws.Range("H141").Value = 1469.8
$ws.Range("I141").Value = 1289.0714
$ws.Range("K141").Value = 3867.2142
$ws.Range("M141").Value = 1312.7858

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12826534
$ws.Range("I32").Value = 12826534
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 12826534
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -12826247
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 5573338.5
$ws.Range("I61").Value = 5573338.5
$ws.Range("K61").Value = 5573338.5
$ws.Range("M61").Value = -5573126.5
$ws.Range("H74").Value = 11890.5
$ws.Range("I74").Value = 14829.667
$ws.Range("J74").Value = 6600
$ws.Range("K74").Value = 14829.667
$ws.Range("L74").Value = 6600
$ws.Range("M74").Value = -13955.667
$ws.Range("N74").Value = -8348
$ws.Range("H77").Value = 11890.5
$ws.Range("I77").Value = 14829.667
$ws.Range("J77").Value = 6600
$ws.Range("K77").Value = 74148.33499999999
$ws.Range("L77").Value = 33000
$ws.Range("M77").Value = -69780.33499999999
$ws.Range("N77").Value = -41736
$ws.Range("H132").Value = 1249689.6
$ws.Range("I132").Value = 1396570.8
$ws.Range("K132").Value = 4189712.4
$ws.Range("M132").Value = -4187182.4
$ws.Range("H136").Value = 5573338.5
$ws.Range("I136").Value = 5573338.5
$ws.Range("K136").Value = 16720015.5
$ws.Range("M136").Value = -16717465.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2912.75
$ws.Range("I86").Value = 3217
$ws.Range("K86").Value = 3217
$ws.Range("M86").Value = -2094
$ws.Range("H89").Value = 2912.75
$ws.Range("I89").Value = 3217
$ws.Range("K89").Value = 16085
$ws.Range("M89").Value = -10469
$ws.Range("H107").Value = 6165.357
$ws.Range("I107").Value = 4291.222
$ws.Range("J107").Value = 9538.799999999999
$ws.Range("K107").Value = 4291.222
$ws.Range("L107").Value = 9538.799999999999
$ws.Range("M107").Value = -2371.222
$ws.Range("N107").Value = -13378.8
$ws.Range("H119").Value = 90761
$ws.Range("J119").Value = 90761
$ws.Range("L119").Value = 90761
$ws.Range("N119").Value = -100437

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 71802400
$ws.Range("I132").Value = 111115620
$ws.Range("J132").Value = 1038624.8
$ws.Range("K132").Value = 333346860
$ws.Range("L132").Value = 3115874.4
$ws.Range("M132").Value = -333344330
$ws.Range("N132").Value = -3120934.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6371.6665
$ws.Range("J68").Value = 7424.1177
$ws.Range("L68").Value = 22272.3531
$ws.Range("N68").Value = -23894.3531
$ws.Range("H71").Value = 6371.6665
$ws.Range("J71").Value = 7424.1177
$ws.Range("L71").Value = 66817.05929999999
$ws.Range("N71").Value = -74929.05929999999
$ws.Range("H75").Value = 7068.75
$ws.Range("J75").Value = 8676.666999999999
$ws.Range("L75").Value = 26030.001
$ws.Range("N75").Value = -28026.001
$ws.Range("H78").Value = 7068.75
$ws.Range("J78").Value = 8676.666999999999
$ws.Range("L78").Value = 78090.003
$ws.Range("N78").Value = -88074.003
$ws.Range("H109").Value = 54896.895
$ws.Range("I109").Value = 73160.07000000001
$ws.Range("J109").Value = 3760
$ws.Range("K109").Value = 219480.21
$ws.Range("L109").Value = 11280
$ws.Range("M109").Value = -218440.21
$ws.Range("N109").Value = -13360
$ws.Range("H110").Value = 12924.909
$ws.Range("I110").Value = 3013.5
$ws.Range("J110").Value = 18588.572
$ws.Range("K110").Value = 9040.5
$ws.Range("L110").Value = 55765.716
$ws.Range("M110").Value = -4950.5
$ws.Range("N110").Value = -63945.716
$ws.Range("H112").Value = 2027
$ws.Range("I112").Value = 2027
$ws.Range("K112").Value = 6081
$ws.Range("M112").Value = -4973
$ws.Range("H116").Value = 2552.8
$ws.Range("I116").Value = 691
$ws.Range("K116").Value = 2073
$ws.Range("M116").Value = 1369
$ws.Range("H121").Value = 20000690
$ws.Range("I121").Value = 33333500
$ws.Range("K121").Value = 100000500
$ws.Range("M121").Value = -99999190
$ws.Range("H129").Value = 1938.2307
$ws.Range("I129").Value = 650.375
$ws.Range("J129").Value = 3998.8
$ws.Range("K129").Value = 1951.125
$ws.Range("L129").Value = 11996.4
$ws.Range("M129").Value = 3048.875
$ws.Range("N129").Value = -21996.4
$ws.Range("H134").Value = 1295.909
$ws.Range("I134").Value = 1295.909
$ws.Range("K134").Value = 3887.727
$ws.Range("M134").Value = 1182.273
$ws.Range("H139").Value = 55829.26
$ws.Range("I139").Value = 65247.25
$ws.Range("J139").Value = 5600
$ws.Range("K139").Value = 195741.75
$ws.Range("L139").Value = 16800
$ws.Range("M139").Value = -190601.75
$ws.Range("N139").Value = -27080

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 3500
$ws.Range("I48").Value = 3500
$ws.Range("K48").Value = 3500
$ws.Range("M48").Value = -3015
$ws.Range("H97").Value = 1168.75
$ws.Range("I97").Value = 926.63635
$ws.Range("K97").Value = 926.63635
$ws.Range("M97").Value = -430.63635
$ws.Range("H107").Value = 49539.855
$ws.Range("I107").Value = 69196
$ws.Range("J107").Value = 399.5
$ws.Range("K107").Value = 69196
$ws.Range("L107").Value = 399.5
$ws.Range("M107").Value = -67276
$ws.Range("N107").Value = -4239.5
$ws.Range("H113").Value = 4499.3335
$ws.Range("I113").Value = 2999.5
$ws.Range("J113").Value = 7499
$ws.Range("K113").Value = 2999.5
$ws.Range("L113").Value = 7499
$ws.Range("M113").Value = -829.5
$ws.Range("N113").Value = -11839

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2573.4285
$ws.Range("I68").Value = 2000.6666
$ws.Range("K68").Value = 2000.6666
$ws.Range("M68").Value = -1251.6666
$ws.Range("H71").Value = 2573.4285
$ws.Range("I71").Value = 2000.6666
$ws.Range("K71").Value = 10003.333
$ws.Range("M71").Value = -6259.333000000001
$ws.Range("H136").Value = 43226.676
$ws.Range("I136").Value = 2655.3
$ws.Range("K136").Value = 7965.900000000001
$ws.Range("M136").Value = -5415.900000000001
$ws.Range("H138").Value = 89000
$ws.Range("J138").Value = 89000
$ws.Range("L138").Value = 89000
$ws.Range("N138").Value = -99280
$ws.Range("H140").Value = 20000
$ws.Range("J140").Value = 20000
$ws.Range("L140").Value = 20000
$ws.Range("N140").Value = -30360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 10949.333
$ws.Range("I8").Value = 6424.5
$ws.Range("J8").Value = 19999
$ws.Range("K8").Value = 6424.5
$ws.Range("L8").Value = 19999
$ws.Range("M8").Value = -6284.5
$ws.Range("N8").Value = -20279
$ws.Range("H62").Value = 2201.6667
$ws.Range("I62").Value = 2301
$ws.Range("J62").Value = 2003
$ws.Range("K62").Value = 2301
$ws.Range("L62").Value = 2003
$ws.Range("M62").Value = -1677
$ws.Range("N62").Value = -3251
$ws.Range("H65").Value = 2201.6667
$ws.Range("I65").Value = 2301
$ws.Range("J65").Value = 2003
$ws.Range("K65").Value = 11505
$ws.Range("L65").Value = 10015
$ws.Range("M65").Value = -8385
$ws.Range("N65").Value = -16255
$ws.Range("H107").Value = 1490.625
$ws.Range("I107").Value = 1048.28
$ws.Range("J107").Value = 3070.4285
$ws.Range("K107").Value = 3144.84
$ws.Range("L107").Value = 9211.2855
$ws.Range("M107").Value = -1224.84
$ws.Range("N107").Value = -13051.2855
$ws.Range("H132").Value = 3874455.2
$ws.Range("I132").Value = 5593053.5
$ws.Range("J132").Value = 7609.1875
$ws.Range("K132").Value = 16779160.5
$ws.Range("L132").Value = 22827.5625
$ws.Range("M132").Value = -16776630.5
$ws.Range("N132").Value = -27887.5625

